$d = $word.ActiveDocument

# Swahili -> English translations of the table labels / cell text.
$replacements = @(
    @("Kichwa cha Video", "Video Title"),
    @("Mada", "Topic"),
    @("Malengo", "Aim(s)"),
    @("Urefu", "Length"),
    @("Mahali pa Kambi", "Camp Location"),
    @("Wawezeshaji", "Facilitators"),
    @("N. ya wanafunzi", "N. of students"),
    @("Tarehe", "Date"),
    @("Rasilimali", "Resources"),
    @("inahitajika", "needed"),
    @("Maandalizi", "Preparations"),
    @("Muda wa video", "Video time"),
    @("Mwezeshaji anafanya nini", "What facilitator does"),
    @("Wanachofanya wanafunzi", "What learners do"),
    @("Utangulizi Mkuu wa Video ya VMC", "General VMC Video Introduction"),
    @("Utangulizi wa Video", "Video Introduction"),
    @("Kusaidia mchakato, kuchochea mawazo", "Assist the process, provoke thoughts"),
    @("Suluhisho", "Solution")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Default document language: Swahili (Kenya) -> Swahili (Tanzania)
$d.Styles("Normal").LanguageID = "sw-TZ"
